$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("G2").Value = 2.15
$ws.Range("I2").Value = 3.3
$ws.Range("J2").Value = 1.03
$ws.Range("K2").Value = 15
$ws.Range("Z2").Value = 15

# Row 7 updates
$ws.Range("H7").Value = 3.3
$ws.Range("I7").Value = 2.5
$ws.Range("L7").Value = 1.33
$ws.Range("M7").Value = 2.8
$ws.Range("N7").Value = 1.98
$ws.Range("O7").Value = 1.65
$ws.Range("Q7").Value = 2.55
$ws.Range("R7").Value = 1.78
$ws.Range("S7").Value = 1.82
$ws.Range("T7").Value = 7.9
$ws.Range("W7").Value = 28
$ws.Range("X7").Value = 23
$ws.Range("Y7").Value = 35
$ws.Range("Z7").Value = 9
$ws.Range("AA7").Value = 6.4
$ws.Range("AB7").Value = 15
$ws.Range("AC7").Value = 75
$ws.Range("AD7").Value = 700
$ws.Range("AE7").Value = 7.8
$ws.Range("AF7").Value = 11.75
$ws.Range("AI7").Value = 22
$ws.Range("AJ7").Value = 35
